$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price (D) and volume-change (E) values from the latest cryptos pull.
# D-column price strings that look numeric (single decimal point) must be forced to
# Text format first, otherwise Excel auto-converts them to numeric cell values.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.371.43'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.317.00'
$ws.Range('E3').Value = '  -2.05%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '312.64'
$ws.Range('E5').Value = '  -4.73%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '106.60'
$ws.Range('E6').Value = '  +6.72%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.629'
$ws.Range('E7').Value = '  -0.95%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.607'
$ws.Range('E9').Value = '  -1.56%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.23'
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0918'
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('E12').Value = '  -1.00%  '
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('E14').Value = '  -2.67%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.60'
$ws.Range('E15').Value = '  -3.97%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.664.23'
$ws.Range('E16').Value = '  -2.21%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.307.28'
$ws.Range('E17').Value = '  -2.56%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '42.278.56'
$ws.Range('E18').Value = '  -0.32%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.48'
$ws.Range('E19').Value = '  -3.58%  '
$ws.Range('E20').Value = '  -0.67%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '75.49'
$ws.Range('E21').Value = '  +0.95%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.51'
$ws.Range('E22').Value = '  -5.83%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '264.22'
$ws.Range('E23').Value = '  -4.46%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.29'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.30'
$ws.Range('E25').Value = '  -2.95%  '
$ws.Range('E26').Value = '  +0.46%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.13'
$ws.Range('E27').Value = '  -2.50%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '23.14'
$ws.Range('E28').Value = '  -2.19%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.24'
$ws.Range('E29').Value = '  +1.53%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.10'
$ws.Range('E30').Value = '  +2.84%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '165.73'
$ws.Range('E31').Value = '  -4.90%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0898'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('E33').Value = '  -6.43%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.91'
$ws.Range('E34').Value = '  -1.27%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.120'
$ws.Range('E35').Value = '  +14.21%  '
$ws.Range('E36').Value = '  -3.59%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.57'
$ws.Range('E37').Value = '  +0.98%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0352'
$ws.Range('E38').Value = '  -1.30%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.72'
$ws.Range('E39').Value = '  -3.15%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.66'
$ws.Range('E40').Value = '  -9.12%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '101.39'
$ws.Range('E41').Value = '  +10.77%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.47'
$ws.Range('E42').Value = '  -3.09%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '70.66'
$ws.Range('E43').Value = '  +2.46%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.233'
$ws.Range('E44').Value = '  +2.33%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.21'
$ws.Range('E46').Value = '  +2.47%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '111.97'
$ws.Range('E47').Value = '  -3.38%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '5.43'
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.05'
$ws.Range('E49').Value = '  -0.41%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '76.00'
$ws.Range('E50').Value = '  +11.14%  '
$ws.Range('E51').Value = '  +0.62%  '
